{"js": "// Replace the date line and every \"A\u00d7B=C\" answer cell in the practice\n// table with its updated value, per the commit's regenerated numbers.\n// Each old value is unique in the document, so a scoped search/replace\n// per pair is safe and keeps formatting (fonts/sizes) untouched since\n// insertText(\"Replace\") only swaps the text of the matched range.\nconst replacements = [\n  [\"2024-03-17 Sunday\", \"2024-03-18 Monday\"],\n  [\"771\u00d74=3084\", \"228\u00d76=1368\"],\n  [\"592\u00d76=3552\", \"924\u00d76=5544\"],\n  [\"525\u00d74=2100\", \"732\u00d74=2928\"],\n  [\"865\u00d72=1730\", \"146\u00d77=1022\"],\n  [\"835\u00d77=5845\", \"370\u00d72=740\"],\n  [\"306\u00d72=612\", \"374\u00d73=1122\"],\n  [\"372\u00d79=3348\", \"476\u00d75=2380\"],\n  [\"976\u00d74=3904\", \"623\u00d73=1869\"],\n  [\"736\u00d72=1472\", \"396\u00d73=1188\"],\n  [\"865\u00d74=3460\", \"510\u00d72=1020\"],\n  [\"221\u00d73=663\", \"221\u00d75=1105\"],\n  [\"878\u00d79=7902\", \"432\u00d72=864\"],\n  [\"354\u00d78=2832\", \"763\u00d77=5341\"],\n  [\"500\u00d75=2500\", \"559\u00d78=4472\"],\n  [\"294\u00d74=1176\", \"267\u00d75=1335\"],\n  [\"395\u00d76=2370\", \"285\u00d79=2565\"],\n  [\"127\u00d73=381\", \"226\u00d76=1356\"],\n  [\"300\u00d75=1500\", \"908\u00d78=7264\"],\n  [\"199\u00d75=995\", \"710\u00d75=3550\"],\n  [\"407\u00d76=2442\", \"474\u00d75=2370\"],\n  [\"233\u00d75=1165\", \"883\u00d73=2649\"],\n  [\"382\u00d74=1528\", \"175\u00d74=700\"],\n  [\"685\u00d73=2055\", \"740\u00d72=1480\"],\n  [\"176\u00d75=880\", \"303\u00d72=606\"],\n  [\"404\u00d77=2828\", \"376\u00d79=3384\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"A\u00d7B=C\" answer cell in the practice\n# table with its updated value, per the commit's regenerated numbers.\n# Each old value is unique in the document, so Find/Replace scoped to\n# the whole story (wdReplaceAll) is safe and leaves formatting intact.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-03-17 Sunday\", \"2024-03-18 Monday\"),\n    @(\"771\u00d74=3084\", \"228\u00d76=1368\"),\n    @(\"592\u00d76=3552\", \"924\u00d76=5544\"),\n    @(\"525\u00d74=2100\", \"732\u00d74=2928\"),\n    @(\"865\u00d72=1730\", \"146\u00d77=1022\"),\n    @(\"835\u00d77=5845\", \"370\u00d72=740\"),\n    @(\"306\u00d72=612\", \"374\u00d73=1122\"),\n    @(\"372\u00d79=3348\", \"476\u00d75=2380\"),\n    @(\"976\u00d74=3904\", \"623\u00d73=1869\"),\n    @(\"736\u00d72=1472\", \"396\u00d73=1188\"),\n    @(\"865\u00d74=3460\", \"510\u00d72=1020\"),\n    @(\"221\u00d73=663\", \"221\u00d75=1105\"),\n    @(\"878\u00d79=7902\", \"432\u00d72=864\"),\n    @(\"354\u00d78=2832\", \"763\u00d77=5341\"),\n    @(\"500\u00d75=2500\", \"559\u00d78=4472\"),\n    @(\"294\u00d74=1176\", \"267\u00d75=1335\"),\n    @(\"395\u00d76=2370\", \"285\u00d79=2565\"),\n    @(\"127\u00d73=381\", \"226\u00d76=1356\"),\n    @(\"300\u00d75=1500\", \"908\u00d78=7264\"),\n    @(\"199\u00d75=995\", \"710\u00d75=3550\"),\n    @(\"407\u00d76=2442\", \"474\u00d75=2370\"),\n    @(\"233\u00d75=1165\", \"883\u00d73=2649\"),\n    @(\"382\u00d74=1528\", \"175\u00d74=700\"),\n    @(\"685\u00d73=2055\", \"740\u00d72=1480\"),\n    @(\"176\u00d75=880\", \"303\u00d72=606\"),\n    @(\"404\u00d77=2828\", \"376\u00d79=3384\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
